$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1208.9231
$ws.Range("I15").Value = 1208.9231
$ws.Range("K15").Value = 3626.7693
$ws.Range("M15").Value = -3457.7693

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 12000.611
$ws.Range("I19").Value = 787.36365
$ws.Range("J19").Value = 29621.428
$ws.Range("K19").Value = 787.36365
$ws.Range("L19").Value = 29621.428
$ws.Range("M19").Value = -612.36365
$ws.Range("N19").Value = -29971.428

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 40003.418
$ws.Range("I32").Value = 63271.145
$ws.Range("J32").Value = 7428.6
$ws.Range("K32").Value = 63271.145
$ws.Range("L32").Value = 7428.6
$ws.Range("M32").Value = -62945.145
$ws.Range("N32").Value = -8080.6

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 25000424
$ws.Range("I39").Value = 439.6154
$ws.Range("J39").Value = 71428970
$ws.Range("K39").Value = 1318.8462
$ws.Range("L39").Value = 214286910
$ws.Range("M39").Value = -1022.8462
$ws.Range("N39").Value = -214287502

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 592.1429000000001
$ws.Range("I98").Value = 513.1389
$ws.Range("J98").Value = 1066.1666
$ws.Range("K98").Value = 513.1389
$ws.Range("L98").Value = 1066.1666
$ws.Range("M98").Value = 984.8611
$ws.Range("N98").Value = -4062.1666

# ALC row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5566.0713
$ws.Range("I100").Value = 6226
$ws.Range("J100").Value = 3146.3333
$ws.Range("K100").Value = 6226
$ws.Range("L100").Value = 3146.3333
$ws.Range("M100").Value = -5685
$ws.Range("N100").Value = -4228.3333

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 13445
$ws.Range("I111").Value = 322.55554
$ws.Range("K111").Value = 967.66662
$ws.Range("M111").Value = 2099.33338

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4222.125
$ws.Range("J112").Value = 3332.6667
$ws.Range("L112").Value = 9998.000100000001
$ws.Range("N112").Value = -12214.0001

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 592.1429000000001
$ws.Range("I122").Value = 513.1389
$ws.Range("J122").Value = 1066.1666
$ws.Range("K122").Value = 1539.4167
$ws.Range("L122").Value = 3198.4998
$ws.Range("M122").Value = 910.5832999999998
$ws.Range("N122").Value = -8098.4998

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3203.5715
$ws.Range("I132").Value = 2791.151
$ws.Range("K132").Value = 8373.453
$ws.Range("M132").Value = -5843.453

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1852.5636
$ws.Range("I138").Value = 1424.975
$ws.Range("J138").Value = 2992.8
$ws.Range("K138").Value = 4274.924999999999
$ws.Range("L138").Value = 8978.400000000001
$ws.Range("M138").Value = 865.0750000000007
$ws.Range("N138").Value = -19258.4

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3521.0286
$ws.Range("I74").Value = 3530.3872
$ws.Range("J74").Value = 3448.5
$ws.Range("K74").Value = 3530.3872
$ws.Range("L74").Value = 3448.5
$ws.Range("M74").Value = -2656.3872
$ws.Range("N74").Value = -5196.5

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3521.0286
$ws.Range("I77").Value = 3530.3872
$ws.Range("J77").Value = 3448.5
$ws.Range("K77").Value = 17651.936
$ws.Range("L77").Value = 17242.5
$ws.Range("M77").Value = -13283.936
$ws.Range("N77").Value = -25978.5

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 40327.45
$ws.Range("I102").Value = 3141.611
$ws.Range("J102").Value = 375000
$ws.Range("K102").Value = 3141.611
$ws.Range("L102").Value = 375000
$ws.Range("M102").Value = -1519.611
$ws.Range("N102").Value = -378244

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 3873.6875
$ws.Range("I110").Value = 4311
$ws.Range("K110").Value = 4311
$ws.Range("M110").Value = -2266

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3028.0908
$ws.Range("I122").Value = 2589.889
$ws.Range("K122").Value = 7769.667
$ws.Range("M122").Value = -5319.667

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1047.164
$ws.Range("I132").Value = 1039.6034
$ws.Range("J132").Value = 1193.3334
$ws.Range("K132").Value = 3118.8102
$ws.Range("L132").Value = 3580.0002
$ws.Range("M132").Value = -588.8101999999999
$ws.Range("N132").Value = -8640.0002

# CRP row 14
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1340

# CRP row 129
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").ClearContents()
$ws.Range("N129").Value = 0

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14395554
$ws.Range("I4").Value = 24977974
$ws.Range("K4").Value = 74933922
$ws.Range("M4").Value = -74933810

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1460.762
$ws.Range("I5").Value = 1010
$ws.Range("K5").Value = 3030
$ws.Range("M5").Value = -2918

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 162.28572
$ws.Range("I17").Value = 164
$ws.Range("K17").Value = 492
$ws.Range("M17").Value = -323

# CUL row 62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 13571.429
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 45000
$ws.Range("N62").Value = -46372

# CUL row 65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 13571.429
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 135000
$ws.Range("N65").Value = -141864

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1446
$ws.Range("I68").Value = 1116.4
$ws.Range("J68").Value = 1610.8
$ws.Range("K68").Value = 3349.2
$ws.Range("L68").Value = 4832.4
$ws.Range("M68").Value = -2538.2
$ws.Range("N68").Value = -6454.4

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1446
$ws.Range("I71").Value = 1116.4
$ws.Range("J71").Value = 1610.8
$ws.Range("K71").Value = 10047.6
$ws.Range("L71").Value = 14497.2
$ws.Range("M71").Value = -5991.6
$ws.Range("N71").Value = -22609.2

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1569.5625
$ws.Range("J113").Value = 1729.8462
$ws.Range("L113").Value = 5189.5386
$ws.Range("N113").Value = -9529.5386

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4449.484
$ws.Range("J131").Value = 5380.6665
$ws.Range("L131").Value = 16141.9995
$ws.Range("N131").Value = -26221.9995

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1460.762
$ws.Range("I135").Value = 1010
$ws.Range("K135").Value = 9090
$ws.Range("M135").Value = -6555

# GSM row 27
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 10000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 10000
$ws.Range("N27").Value = -10332

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1785.1428
$ws.Range("I102").Value = 1604.6842
$ws.Range("K102").Value = 1604.6842
$ws.Range("M102").Value = 17.31580000000008

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3468.5095
$ws.Range("I132").Value = 2874.3062
$ws.Range("K132").Value = 8622.918600000001
$ws.Range("M132").Value = -6092.918600000001

# GSM row 133
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 60879.6
$ws.Range("J133").Value = 61225
$ws.Range("L133").Value = 61225
$ws.Range("N133").Value = -71345

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1818.091
$ws.Range("I2").Value = 1818.091
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1818.091
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -1706.091

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2936486.2
$ws.Range("I22").Value = 418.2
$ws.Range("J22").Value = 5383209.5
$ws.Range("K22").Value = 418.2
$ws.Range("L22").Value = 5383209.5
$ws.Range("M22").Value = -123.2
$ws.Range("N22").Value = -5383799.5

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2936486.2
$ws.Range("I27").Value = 418.2
$ws.Range("J27").Value = 5383209.5
$ws.Range("K27").Value = 418.2
$ws.Range("L27").Value = 5383209.5
$ws.Range("M27").Value = -311.2
$ws.Range("N27").Value = -5383423.5

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2346.4
$ws.Range("I61").Value = 2346.4
$ws.Range("K61").Value = 2346.4
$ws.Range("M61").Value = -2144.4

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3232.6667
$ws.Range("I93").Value = 3232.6667
$ws.Range("K93").Value = 3232.6667
$ws.Range("M93").Value = -1984.6667

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2346.4
$ws.Range("I113").Value = 2346.4
$ws.Range("K113").Value = 2346.4
$ws.Range("M113").Value = -176.4000000000001

# LTW row 121
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 27000
$ws.Range("J121").Value = 27000
$ws.Range("L121").Value = 27000
$ws.Range("N121").Value = -30494

# LTW row 125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 200714.5
$ws.Range("J125").Value = 200714.5
$ws.Range("L125").Value = 200714.5
$ws.Range("N125").Value = -210554.5

# WVR row 2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2508335
$ws.Range("I2").Value = 2508335
$ws.Range("K2").Value = 2508335
$ws.Range("M2").Value = -2508223

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 90000
$ws.Range("J46").Value = 90000
$ws.Range("L46").Value = 90000
$ws.Range("N46").Value = -90462

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5868.75
$ws.Range("I96").Value = 2630.4
$ws.Range("J96").Value = 11266
$ws.Range("K96").Value = 2630.4
$ws.Range("L96").Value = 11266
$ws.Range("M96").Value = -1257.4
$ws.Range("N96").Value = -14012

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2591.8572
$ws.Range("I113").Value = 434.9091
$ws.Range("K113").Value = 1304.7273
$ws.Range("M113").Value = 865.2727

# WVR row 123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49999.168
$ws.Range("J123").Value = 49999.168
$ws.Range("L123").Value = 49999.168
$ws.Range("N123").Value = -59799.168

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1142.8334
$ws.Range("I126").Value = 981.96
$ws.Range("K126").Value = 2945.88
$ws.Range("M126").Value = -475.8800000000001

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6106.706
$ws.Range("I132").Value = 6219.8184
$ws.Range("K132").Value = 18659.4552
$ws.Range("M132").Value = -16129.4552

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 270000
$ws.Range("N134").Value = -275070
